$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.95419801644062
$ws.Range("C2").Value = 13.16508314515939
$ws.Range("E2").Value = 17.10419056647336
$ws.Range("F2").Value = 39.22751005504706
$ws.Range("G2").Value = 34.25011559472365
$ws.Range("H2").Value = 15.83807710009416
$ws.Range("J2").Value = 8.407879324138483
$ws.Range("K2").Value = 7.741438133827988
$ws.Range("L2").Value = 12.25328226195592
$ws.Range("M2").Value = 15.31431653136852
$ws.Range("O2").Value = 24.80297989727734
$ws.Range("B3").Value = 11.71148820215922
$ws.Range("C3").Value = 13.20282452466995
$ws.Range("E3").Value = 17.12567554219492
$ws.Range("F3").Value = 39.29857082279469
$ws.Range("G3").Value = 34.40081413437901
$ws.Range("H3").Value = 15.89102589592942
$ws.Range("J3").Value = 8.394822284307496
$ws.Range("K3").Value = 7.497271212481664
$ws.Range("L3").Value = 12.24678709144198
$ws.Range("M3").Value = 15.26491782407435
$ws.Range("O3").Value = 24.89979808133203
$ws.Range("B4").Value = 11.56142743859852
$ws.Range("C4").Value = 13.22722331977866
$ws.Range("E4").Value = 17.14118354540987
$ws.Range("F4").Value = 39.35035707735191
$ws.Range("G4").Value = 34.50185722361358
$ws.Range("H4").Value = 15.92564285619333
$ws.Range("J4").Value = 8.386726299365327
$ws.Range("K4").Value = 7.341927430981277
$ws.Range("L4").Value = 12.24423233585893
$ws.Range("M4").Value = 15.23629770769538
$ws.Range("O4").Value = 24.96355512581524
$ws.Range("B5").Value = 11.50009715062767
$ws.Range("C5").Value = 13.2374750353853
$ws.Range("E5").Value = 17.14808623472776
$ws.Range("F5").Value = 39.3735090615497
$ws.Range("G5").Value = 34.54516858944622
$ws.Range("H5").Value = 15.94027979561418
$ws.Range("J5").Value = 8.383407687479412
$ws.Range("K5").Value = 7.2773199003947
$ws.Range("L5").Value = 12.24355339227092
$ws.Range("M5").Value = 15.22507254280053
$ws.Range("O5").Value = 24.99062042996326
$ws.Range("B6").Value = 11.48990489148325
$ws.Range("C6").Value = 13.23919601373018
$ws.Range("E6").Value = 17.14926765309481
$ws.Range("F6").Value = 39.37747710058054
$ws.Range("G6").Value = 34.55248923001365
$ws.Range("H6").Value = 15.94274229257914
$ws.Range("J6").Value = 8.382855472721843
$ws.Range("K6").Value = 7.266515007897604
$ws.Range("L6").Value = 12.24346257448515
$ws.Range("M6").Value = 15.22323526588077
$ws.Range("O6").Value = 24.99518005992126
$ws.Range("B7").Value = 11.56060093390504
$ws.Range("C7").Value = 13.22736032568601
$ws.Range("E7").Value = 17.14127427594385
$ws.Range("F7").Value = 39.35066102099126
$ws.Range("G7").Value = 34.50243269642363
$ws.Range("H7").Value = 15.92583810728534
$ws.Range("J7").Value = 8.386681621549757
$ws.Range("K7").Value = 7.341061305789437
$ws.Range("L7").Value = 12.24422171091004
$ws.Range("M7").Value = 15.23614453900811
$ws.Range("O7").Value = 24.96391574981761
$ws.Range("B8").Value = 11.87077463220537
$ws.Range("C8").Value = 13.17784255388623
$ws.Range("E8").Value = 17.11111827623913
$ws.Range("F8").Value = 39.25031831032486
$ws.Range("G8").Value = 34.30030628388381
$ws.Range("H8").Value = 15.85589708453018
$ws.Range("J8").Value = 8.403393879651588
$ws.Range("K8").Value = 7.658409569086755
$ws.Range("L8").Value = 12.25074631294165
$ws.Range("M8").Value = 15.29693339619251
$ws.Range("O8").Value = 24.83546816659543
$ws.Range("B9").Value = 12.46737197212512
$ws.Range("C9").Value = 13.09042214627022
$ws.Range("E9").Value = 17.07032994211212
$ws.Range("F9").Value = 39.11831750096734
$ws.Range("G9").Value = 33.97172311161781
$ws.Range("H9").Value = 15.73542603403005
$ws.Range("J9").Value = 8.43552328387384
$ws.Range("K9").Value = 8.23533307753749
$ws.Range("L9").Value = 12.27483414319606
$ws.Range("M9").Value = 15.42936953201598
$ws.Range("O9").Value = 24.61778580216651
$ws.Range("B10").Value = 12.89421546117151
$ws.Range("C10").Value = 13.03204243048264
$ws.Range("E10").Value = 17.05150449763593
$ws.Range("F10").Value = 39.06088904684616
$ws.Range("G10").Value = 33.77195333416126
$ws.Range("H10").Value = 15.6570467077664
$ws.Range("J10").Value = 8.45871480254366
$ws.Range("K10").Value = 8.628758187306559
$ws.Range("L10").Value = 12.29930376245124
$ws.Range("M10").Value = 15.53425228306412
$ws.Range("O10").Value = 24.47870909332185
$ws.Range("B11").Value = 13.08505278738018
$ws.Range("C11").Value = 13.00674207637681
$ws.Range("E11").Value = 17.04534870292972
$ws.Range("F11").Value = 39.0433551550165
$ws.Range("G11").Value = 33.69018069537982
$ws.Range("H11").Value = 15.62358151694376
$ws.Range("J11").Value = 8.469170651977178
$ws.Range("K11").Value = 8.800645067880811
$ws.Range("L11").Value = 12.31187966227449
$ws.Range("M11").Value = 15.58350293865277
$ws.Range("O11").Value = 24.41996914041324
$ws.Range("B12").Value = 13.1567725590379
$ws.Range("C12").Value = 12.99734134910327
$ws.Range("E12").Value = 17.04336283125207
$ws.Range("F12").Value = 39.03795006233908
$ws.Range("G12").Value = 33.66052997759339
$ws.Range("H12").Value = 15.61122346352299
$ws.Range("J12").Value = 8.4731160337943
$ws.Range("K12").Value = 8.86468126281361
$ws.Range("L12").Value = 12.31684705762377
$ws.Range("M12").Value = 15.60236440698529
$ws.Range("O12").Value = 24.39837705092885
$ws.Range("B13").Value = 13.14135172764315
$ws.Range("C13").Value = 12.99935797196207
$ws.Range("E13").Value = 17.04377518977858
$ws.Range("F13").Value = 39.03905925378141
$ws.Range("G13").Value = 33.66685722424824
$ws.Range("H13").Value = 15.61387101111098
$ws.Range("J13").Value = 8.472266956730195
$ws.Range("K13").Value = 8.850937295041106
$ws.Range("L13").Value = 12.31576815693081
$ws.Range("M13").Value = 15.59829301590527
$ws.Range("O13").Value = 24.40299831161039
$ws.Range("B14").Value = 13.09096452806346
$ws.Range("C14").Value = 13.00596507043789
$ws.Range("E14").Value = 17.04517841236034
$ws.Range("F14").Value = 39.04288573753414
$ws.Range("G14").Value = 33.68771493004773
$ws.Range("H14").Value = 15.62255851202416
$ws.Range("J14").Value = 8.469495524763065
$ws.Range("K14").Value = 8.805934665825372
$ws.Range("L14").Value = 12.31228423593431
$ws.Range("M14").Value = 15.5850505061092
$ws.Range("O14").Value = 24.41817968703778
$ws.Range("B15").Value = 13.06002788422574
$ws.Range("C15").Value = 13.01003552347099
$ws.Range("E15").Value = 17.04608284766494
$ws.Range("F15").Value = 39.04539032089191
$ws.Range("G15").Value = 33.70066227101123
$ws.Range("H15").Value = 15.62792080044572
$ws.Range("J15").Value = 8.467796101955114
$ws.Range("K15").Value = 8.778231076392697
$ws.Range("L15").Value = 12.31017687567316
$ws.Range("M15").Value = 15.57696630483244
$ws.Range("O15").Value = 24.42756357230936
$ws.Range("B16").Value = 12.88167088567293
$ws.Range("C16").Value = 13.03372109367157
$ws.Range("E16").Value = 17.05195515887619
$ws.Range("F16").Value = 39.06220774619179
$ws.Range("G16").Value = 33.77748104820921
$ws.Range("H16").Value = 15.65927776754606
$ws.Range("J16").Value = 8.458029553010739
$ws.Range("K16").Value = 8.617379339690237
$ws.Range("L16").Value = 12.2985107509142
$ws.Range("M16").Value = 15.53106375725596
$ws.Range("O16").Value = 24.48263902059697
$ws.Range("B17").Value = 12.77135158569526
$ws.Range("C17").Value = 13.04857278895894
$ws.Range("E17").Value = 17.0561736449398
$ws.Range("F17").Value = 39.07472463978224
$ws.Range("G17").Value = 33.82694264496436
$ws.Range("H17").Value = 15.67907484460671
$ws.Range("J17").Value = 8.452013645464532
$ws.Range("K17").Value = 8.516861358588908
$ws.Range("L17").Value = 12.29172217274064
$ws.Range("M17").Value = 15.50329117820639
$ws.Range("O17").Value = 24.51758579126104
$ws.Range("B18").Value = 12.70758855850255
$ws.Range("C18").Value = 13.05723343938724
$ws.Range("E18").Value = 17.05882672875884
$ws.Range("F18").Value = 39.08273269411263
$ws.Range("G18").Value = 33.85624819035531
$ws.Range("H18").Value = 15.69066775336687
$ws.Range("J18").Value = 8.448544628290762
$ws.Range("K18").Value = 8.458381083903809
$ws.Range("L18").Value = 12.28795369070168
$ws.Range("M18").Value = 15.48746259569888
$ws.Range("O18").Value = 24.53811229531043
$ws.Range("B19").Value = 12.68594829694164
$ws.Range("C19").Value = 13.06018613678222
$ws.Range("E19").Value = 17.05976398817312
$ws.Range("F19").Value = 39.08558298488696
$ws.Range("G19").Value = 33.86631748043204
$ws.Range("H19").Value = 15.69462833337768
$ws.Range("J19").Value = 8.447368577636219
$ws.Range("K19").Value = 8.438467589646596
$ws.Range("L19").Value = 12.28670120315388
$ws.Range("M19").Value = 15.48212859416351
$ws.Range("O19").Value = 24.54513537531389
$ws.Range("B20").Value = 12.78312787180036
$ws.Range("C20").Value = 13.04697955765414
$ws.Range("E20").Value = 17.05570112245866
$ws.Range("F20").Value = 39.07330850416231
$ws.Range("G20").Value = 33.82158868274146
$ws.Range("H20").Value = 15.67694607895223
$ws.Range("J20").Value = 8.452654966117729
$ws.Range("K20").Value = 8.52763075983626
$ws.Range("L20").Value = 12.29243075772874
$ws.Range("M20").Value = 15.50623263250712
$ws.Range("O20").Value = 24.51382155003773
$ws.Range("B21").Value = 13.10577980183668
$ws.Range("C21").Value = 13.00401952726929
$ws.Range("E21").Value = 17.04475689244029
$ws.Range("F21").Value = 39.04172830813239
$ws.Range("G21").Value = 33.68155278790904
$ws.Range("H21").Value = 15.61999825061599
$ws.Range("J21").Value = 8.470309945998828
$ws.Range("K21").Value = 8.819181881925235
$ws.Range("L21").Value = 12.31330200055811
$ws.Range("M21").Value = 15.58893449815189
$ws.Range("O21").Value = 24.41370286452751
$ws.Range("B22").Value = 13.31344011232527
$ws.Range("C22").Value = 12.97699130298645
$ws.Range("E22").Value = 17.03961569368548
$ws.Range("F22").Value = 39.02828432765065
$ws.Range("G22").Value = 33.59769748673489
$ws.Range("H22").Value = 15.58461244321329
$ws.Range("J22").Value = 8.48176653575872
$ws.Range("K22").Value = 9.003572680122042
$ws.Range("L22").Value = 12.32813720655416
$ws.Range("M22").Value = 15.64421183287367
$ws.Range("O22").Value = 24.35206695373118
$ws.Range("B23").Value = 13.20292259411929
$ws.Range("C23").Value = 12.99132107621767
$ws.Range("E23").Value = 17.04217597683464
$ws.Range("F23").Value = 39.03480163459518
$ws.Range("G23").Value = 33.64174932293255
$ws.Range("H23").Value = 15.60333094156377
$ws.Range("J23").Value = 8.475659601598609
$ws.Range("K23").Value = 8.905733471142685
$ws.Range("L23").Value = 12.32011094404485
$ws.Range("M23").Value = 15.61460038585163
$ws.Range("O23").Value = 24.38461555742899
$ws.Range("B24").Value = 12.777804863639
$ws.Range("C24").Value = 13.04769947678798
$ws.Range("E24").Value = 17.05591404016304
$ws.Range("F24").Value = 39.07394620995868
$ws.Range("G24").Value = 33.82400650132723
$ws.Range("H24").Value = 15.67790783566557
$ws.Range("J24").Value = 8.45236505738035
$ws.Range("K24").Value = 8.522764063632913
$ws.Range("L24").Value = 12.29210998774507
$ws.Range("M24").Value = 15.50490236965167
$ws.Range("O24").Value = 24.51552200809545
$ws.Range("B25").Value = 12.30768932978073
$ws.Range("C25").Value = 13.11304072901784
$ws.Range("E25").Value = 17.07940416242964
$ws.Range("F25").Value = 39.14708402254485
$ws.Range("G25").Value = 34.0533236079984
$ws.Range("H25").Value = 15.76623476699296
$ws.Range("J25").Value = 8.42690251912715
$ws.Range("K25").Value = 8.084425830444244
$ws.Range("L25").Value = 12.27483414319606
$ws.Range("M25").Value = 15.39217426300592
$ws.Range("O25").Value = 24.67301271998919
